$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (hunk 0)
$ws.Range("H19").Value = 2766.5557
$ws.Range("I19").Value = 2900.3333
$ws.Range("J19").Value = 2499
$ws.Range("K19").Value = 2900.3333
$ws.Range("L19").Value = 2499
$ws.Range("M19").Value = -2725.3333
$ws.Range("N19").Value = -2849

# Row 51 (hunk 1)
$ws.Range("H51").Value = 4998
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 4997
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 4997
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -5965

# Row 103 (hunk 2)
$ws.Range("H103").Value = 1466.3334
$ws.Range("I103").Value = 1500
$ws.Range("J103").Value = 1399
$ws.Range("K103").Value = 4500
$ws.Range("L103").Value = 4197
$ws.Range("M103").Value = -3914
$ws.Range("N103").Value = -5369

# Row 132 (hunk 3)
$ws.Range("H132").Value = 1062.9584
$ws.Range("I132").Value = 1062.9584
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3188.8752
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -658.8751999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (hunk 4)
$ws.Range("H2").Value = 2197.8
$ws.Range("I2").Value = 2197.8
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2197.8
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2084.8

# Row 4 (hunk 5)
$ws.Range("H4").Value = 2995.6667
$ws.Range("I4").Value = 1994
$ws.Range("J4").Value = 3496.5
$ws.Range("K4").Value = 1994
$ws.Range("L4").Value = 3496.5
$ws.Range("M4").Value = -1878
$ws.Range("N4").Value = -3728.5

# Row 32 (hunk 6)
$ws.Range("H32").Value = 5953.46
$ws.Range("I32").Value = 6056.592
$ws.Range("J32").Value = 900
$ws.Range("K32").Value = 6056.592
$ws.Range("L32").Value = 900
$ws.Range("M32").Value = -5769.592
$ws.Range("N32").Value = -1474

# Row 63 (hunk 7)
$ws.Range("H63").Value = 1333.3334
$ws.Range("I63").Value = 1333.3334
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1333.3334
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -647.3334

# Row 66 (hunk 8)
$ws.Range("H66").Value = 1333.3334
$ws.Range("I66").Value = 1333.3334
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 6666.666999999999
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -3234.666999999999

# Row 97 (hunk 9)
$ws.Range("H97").Value = 312.18182
$ws.Range("I97").Value = 243.4
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 243.4
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 252.6
$ws.Range("N97").Value = -1992

# Row 116 (hunk 10)
$ws.Range("H116").Value = 2197.8
$ws.Range("I116").Value = 2197.8
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2197.8
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 96.19999999999982

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (hunk 11)
$ws.Range("H3").Value = 2197.8
$ws.Range("I3").Value = 2197.8
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2197.8
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -2083.8

# Row 86 (hunk 12)
$ws.Range("H86").Value = 2490.3845
$ws.Range("I86").Value = 2281.25
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2281.25
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1158.25
$ws.Range("N86").Value = -7246

# Row 89 (hunk 13)
$ws.Range("H89").Value = 2490.3845
$ws.Range("I89").Value = 2281.25
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 11406.25
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -5790.25
$ws.Range("N89").Value = -36232

# Row 105 (hunk 14)
$ws.Range("H105").Value = 5591.7144
$ws.Range("I105").Value = 4995.2
$ws.Range("J105").Value = 7083
$ws.Range("K105").Value = 4995.2
$ws.Range("L105").Value = 7083
$ws.Range("M105").Value = -3248.2
$ws.Range("N105").Value = -10577

# Row 134 (hunk 15)
$ws.Range("H134").Value = 2761.2144
$ws.Range("I134").Value = 2715.7
$ws.Range("J134").Value = 2875
$ws.Range("K134").Value = 8147.099999999999
$ws.Range("L134").Value = 8625
$ws.Range("M134").Value = -5612.099999999999
$ws.Range("N134").Value = -13695

$ws = $wb.Worksheets.Item("CRP")
# Row 132 (hunk 16)
$ws.Range("H132").Value = 4997.3335
$ws.Range("I132").Value = 4996.5
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 14989.5
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -12459.5
$ws.Range("N132").Value = -20057

# Row 134 (hunk 17)
$ws.Range("H134").Value = 3642.1428
$ws.Range("I134").Value = 3831.6667
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 11495.0001
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -8960.000100000001
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (hunk 18)
$ws.Range("H5").Value = 478
$ws.Range("I5").Value = 424.6
$ws.Range("J5").Value = 745
$ws.Range("K5").Value = 1273.8
$ws.Range("L5").Value = 2235
$ws.Range("M5").Value = -1161.8
$ws.Range("N5").Value = -2459

# Row 8 (hunk 19)
$ws.Range("H8").Value = 499.66666
$ws.Range("I8").Value = 499.66666
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1498.99998
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1359.99998

# Row 107 (hunk 20)
$ws.Range("H107").Value = 237
$ws.Range("I107").Value = 3
$ws.Range("J107").Value = 283.8
$ws.Range("K107").Value = 9
$ws.Range("L107").Value = 851.4000000000001
$ws.Range("M107").Value = 1911
$ws.Range("N107").Value = -4691.4

# Row 135 (hunk 21)
$ws.Range("H135").Value = 478
$ws.Range("I135").Value = 424.6
$ws.Range("J135").Value = 745
$ws.Range("K135").Value = 3821.4
$ws.Range("L135").Value = 6705
$ws.Range("M135").Value = -1286.4
$ws.Range("N135").Value = -11775

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (hunk 22)
$ws.Range("H80").Value = 10000.5
$ws.Range("I80").Value = 9998.666999999999
$ws.Range("J80").Value = 10006
$ws.Range("K80").Value = 9998.666999999999
$ws.Range("L80").Value = 10006
$ws.Range("M80").Value = -9000.666999999999
$ws.Range("N80").Value = -12002

# Row 83 (hunk 23)
$ws.Range("H83").Value = 10000.5
$ws.Range("I83").Value = 9998.666999999999
$ws.Range("J83").Value = 10006
$ws.Range("K83").Value = 49993.335
$ws.Range("L83").Value = 50030
$ws.Range("M83").Value = -45001.335
$ws.Range("N83").Value = -60014

# Row 126 (hunk 24)
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0

# Row 132 (hunk 25)
$ws.Range("H132").Value = 4499.5713
$ws.Range("I132").Value = 4500
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -18558.5

# Row 134 (hunk 26)
$ws.Range("H134").Value = 61137.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 61137.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 183412.5
$ws.Range("N134").Value = -188482.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (hunk 27)
$ws.Range("H22").Value = 4574.9
$ws.Range("I22").Value = 1616.3334
$ws.Range("J22").Value = 5842.857
$ws.Range("K22").Value = 1616.3334
$ws.Range("L22").Value = 5842.857
$ws.Range("M22").Value = -1321.3334
$ws.Range("N22").Value = -6432.857

# Row 27 (hunk 28)
$ws.Range("H27").Value = 4574.9
$ws.Range("I27").Value = 1616.3334
$ws.Range("J27").Value = 5842.857
$ws.Range("K27").Value = 1616.3334
$ws.Range("L27").Value = 5842.857
$ws.Range("M27").Value = -1509.3334
$ws.Range("N27").Value = -6056.857

# Row 55 (hunk 29)
$ws.Range("H55").Value = 458
$ws.Range("I55").Value = 410.66666
$ws.Range("J55").Value = 600
$ws.Range("K55").Value = 410.66666
$ws.Range("L55").Value = 600
$ws.Range("M55").Value = -237.66666
$ws.Range("N55").Value = -946

# Row 132 (hunk 30)
$ws.Range("H132").Value = 5999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
# Row 54 (hunk 31)
$ws.Range("H54").Value = 9099.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 9099.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 9099.5
$ws.Range("N54").Value = -10139.5

# Row 122 (hunk 32)
$ws.Range("H122").Value = 3899
$ws.Range("I122").Value = 3899
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11697
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9247

# Row 132 (hunk 33)
$ws.Range("H132").Value = 4798
$ws.Range("I132").Value = 4597
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 13791
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -11261
$ws.Range("N132").Value = -20057

# Row 136 (hunk 34)
$ws.Range("H136").Value = 3286.25
$ws.Range("I136").Value = 3162.5
$ws.Range("J136").Value = 3471.875
$ws.Range("K136").Value = 9487.5
$ws.Range("L136").Value = 10415.625
$ws.Range("M136").Value = -6937.5
$ws.Range("N136").Value = -15515.625
